# Update recomputed NATMI TPM-derived statistics for the Cp-Slc40a1 ligand-
# receptor pair sheet (new TPM values change ligand/receptor expression
# and all downstream specificity/edge-weight figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 4.694776333333333
$ws.Range("H2").Value = 14.084329
$ws.Range("I2").Value = 0.03090644368652182
$ws.Range("J2").Value = 0.03090644368652182
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03428266666666667
$ws.Range("N2").Value = 0.102848
$ws.Range("O2").Value = 0.01042589467522935
$ws.Range("P2").Value = 0.01042589467522935
$ws.Range("Q2").Value = 0.1609494521102222
$ws.Range("R2").Value = 1.448545068992
$ws.Range("S2").Value = 0.0003222273266615836
$ws.Range("T2").Value = 0.0003222273266615836
$ws.Range("G3").Value = 4.694776333333333
$ws.Range("H3").Value = 14.084329
$ws.Range("I3").Value = 0.03090644368652182
$ws.Range("J3").Value = 0.03090644368652182
$ws.Range("O3").Value = 0.5095441114141793
$ws.Range("P3").Value = 0.5095441114141792
$ws.Range("Q3").Value = 7.866072707692888
$ws.Range("R3").Value = 70.794654369236
$ws.Range("S3").Value = 0.01574819638522114
$ws.Range("T3").Value = 0.01574819638522113
$ws.Range("G4").Value = 4.694776333333333
$ws.Range("H4").Value = 14.084329
$ws.Range("I4").Value = 0.03090644368652182
$ws.Range("J4").Value = 0.03090644368652182
$ws.Range("M4").Value = 1.578445666666666
$ws.Range("N4").Value = 4.735336999999999
$ws.Range("O4").Value = 0.4800299939105914
$ws.Range("P4").Value = 0.4800299939105914
$ws.Range("Q4").Value = 7.410449359319221
$ws.Range("R4").Value = 66.694044233873
$ws.Range("S4").Value = 0.01483601997463911
$ws.Range("T4").Value = 0.01483601997463911
$ws.Range("I5").Value = 0.609881555489188
$ws.Range("J5").Value = 0.609881555489188
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.03428266666666667
$ws.Range("N5").Value = 0.102848
$ws.Range("O5").Value = 0.01042589467522935
$ws.Range("P5").Value = 0.01042589467522935
$ws.Range("Q5").Value = 3.176040025948445
$ws.Range("R5").Value = 28.584360233536
$ws.Range("S5").Value = 0.006358560861895318
$ws.Range("T5").Value = 0.006358560861895317
$ws.Range("I6").Value = 0.609881555489188
$ws.Range("J6").Value = 0.609881555489188
$ws.Range("O6").Value = 0.5095441114141793
$ws.Range("P6").Value = 0.5095441114141792
$ws.Range("S6").Value = 0.3107615552596358
$ws.Range("T6").Value = 0.3107615552596357
$ws.Range("I7").Value = 0.609881555489188
$ws.Range("J7").Value = 0.609881555489188
$ws.Range("M7").Value = 1.578445666666666
$ws.Range("N7").Value = 4.735336999999999
$ws.Range("O7").Value = 0.4800299939105914
$ws.Range("P7").Value = 0.4800299939105914
$ws.Range("Q7").Value = 146.2315246611954
$ws.Range("S7").Value = 0.2927614393676569
$ws.Range("T7").Value = 0.2927614393676569
$ws.Range("G8").Value = 54.56532033333334
$ws.Range("I8").Value = 0.3592120008242901
$ws.Range("J8").Value = 0.3592120008242901
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.03428266666666667
$ws.Range("N8").Value = 0.102848
$ws.Range("O8").Value = 0.01042589467522935
$ws.Range("P8").Value = 0.01042589467522935
$ws.Range("Q8").Value = 1.870644688547556
$ws.Range("R8").Value = 16.835802196928
$ws.Range("S8").Value = 0.003745106486672447
$ws.Range("T8").Value = 0.003745106486672446
$ws.Range("G9").Value = 54.56532033333334
$ws.Range("I9").Value = 0.3592120008242901
$ws.Range("J9").Value = 0.3592120008242901
$ws.Range("O9").Value = 0.5095441114141793
$ws.Range("P9").Value = 0.5095441114141792
$ws.Range("Q9").Value = 91.42390320345822
$ws.Range("R9").Value = 822.8151288311241
$ws.Range("S9").Value = 0.1830343597693224
$ws.Range("T9").Value = 0.1830343597693223
$ws.Range("G10").Value = 54.56532033333334
$ws.Range("I10").Value = 0.3592120008242901
$ws.Range("J10").Value = 0.3592120008242901
$ws.Range("M10").Value = 1.578445666666666
$ws.Range("N10").Value = 4.735336999999999
$ws.Range("O10").Value = 0.4800299939105914
$ws.Range("P10").Value = 0.4800299939105914
$ws.Range("Q10").Value = 86.12839343042855
$ws.Range("S10").Value = 0.1724325345682954
$ws.Range("T10").Value = 0.1724325345682954
